# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" worksheet (cloned from "2022-Q1" so the text
#    formatting of the numeric-looking fund figures is preserved), placed
#    right after the "总计" summary sheet.
# 2. Insert a new row into the "总计" summary sheet for the 2022-Q4 totals.

$wb = $excel.ActiveWorkbook

# --- 1. Build the new "2022-Q4" sheet from a copy of "2022-Q1" ---------
$template = $wb.Worksheets.Item("2022-Q1")
$template.Copy($template)
$newSheet = $wb.Worksheets.Item("2022-Q1 (2)")
$newSheet.Name = "2022-Q4"

# Only one fund is reported for 2022-Q4, so drop the second data row.
$newSheet.Range("A3").EntireRow.Delete()

# Fund code / name / rank are unchanged from the 2022-Q1 snapshot; only the
# size / position figures differ. Force the target cells to text first (the
# source cells were text too) so the numeric-looking strings round-trip as
# text instead of being coerced to numbers, then strip the stray
# NumberFormat that introduces by repainting the format from a plain cell.
$newSheet.Range("D2:G2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.25"
$newSheet.Range("E2").Value = "99.41"
$newSheet.Range("F2").Value = "2.34"
$newSheet.Range("G2").Value = "0.0058"
$newSheet.Range("H2").Copy()
$newSheet.Range("D2:G2").PasteSpecial(-4122)

# --- 2. Add the 2022-Q4 row to the "总计" summary sheet -----------------
$summary = $wb.Worksheets.Item("总计")
$summary.Range("A2").EntireRow.Insert()

# Re-use the index column's formatting (border + centered bold) from the
# row below, then wipe the inherited format off the rest of the new row.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.01

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q1"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.02

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2021-Q4"
$summary.Range("C4").Value = 1
$summary.Range("D4").Value = 0.04

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q2"
$summary.Range("C5").Value = 1
$summary.Range("D5").Value = 0.48

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q1"
$summary.Range("C6").Value = 1
$summary.Range("D6").Value = 0.04

# Keep the originally-selected sheet ("2021-Q1", last tab) active instead of
# leaving the freshly-created/copied sheet selected.
$wb.Worksheets.Item("2021-Q1").Activate()
